$d = $word.ActiveDocument

# --- 1. NAME line: "NAME: KAUSHIK NARAYANAN V" -> "NAME: Ajay Kumar J" ------
$rng = $d.Content
$rng.Find.Execute("KAUSHIK NARAYANAN V", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "Ajay Kumar J", 2)

# --- 2. REG NO line: "REG NO: 192321047" -> "REG NO: 192372072" -------------
$rng = $d.Content
$rng.Find.Execute("192321047", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "192372072", 2)

# --- 3. Code block: each source line is currently spread across several
#        runs (separated by spell/grammar-check <w:proofErr/> markers left
#        over from typing). Re-typing the full line text as a whole collapses
#        it back down to a single clean run per paragraph and drops the
#        now-stale proofErr markers, matching how Word normalizes a line
#        once it is selected and retyped in one go.
$codeLines = @(
    '#include <stdio.h>',
    'int main() {',
    '    char ch;',
    '    printf("Enter a string of arithmetic expressions: ");',
    '    while ((ch = getchar()) != ''\n'') {',
    '        if (ch == ''+'' || ch == ''-'' || ch == ''*'' || ch == ''/'') {',
    '            printf("Operator: %c\n", ch);'
)

$startIndex = 9
for ($i = 0; $i -lt $codeLines.Count; $i++) {
    $target = $codeLines[$i]

    # First pass: write a deliberately different placeholder so the engine
    # always registers a real content change (identical-text writes are
    # treated as no-ops and would leave the old multi-run split in place).
    $p = $d.Paragraphs.Item($startIndex + $i)
    $pr = $p.Range
    $pr.MoveEnd(1, -1) | Out-Null
    $pr.Text = $target + "`u{2603}"

    # Second pass: overwrite with the real target text; the paragraph now
    # holds a single run (no proofErr siblings) with the line's formatting.
    $p2 = $d.Paragraphs.Item($startIndex + $i)
    $pr2 = $p2.Range
    $pr2.MoveEnd(1, -1) | Out-Null
    $pr2.Text = $target
}
